$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.959.39'
$ws.Range('E2').Value = '  +1.36%  '

$ws.Range('D3').Value = '2.571.85'
$ws.Range('E3').Value = '  +2.52%  '

$ws.Range('E4').Value = '  +0.31%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.49'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.60%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.95'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.83%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.576'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.58%  '

$ws.Range('E8').Value = '  +0.16%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.549'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.19%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.41'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.55%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0808'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.62%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.64'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.58%  '

$ws.Range('E13').Value = '  +7.74%  '

$ws.Range('D14').Value = '2.579.60'
$ws.Range('E14').Value = '  +3.60%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.884'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.89%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.38'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.23%  '

$ws.Range('D17').Value = '42.953.94'
$ws.Range('E17').Value = '  +1.60%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.12'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +7.33%  '

$ws.Range('D19').Value = '0.0₃0992'
$ws.Range('E19').Value = '  +3.90%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.66'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.23%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.31%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '254.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.68%  '

$ws.Range('E23').Value = '  +3.99%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.12'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.35%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '28.58'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.37%  '

$ws.Range('E26').Value = '  -0.06%  '

$ws.Range('E27').Value = '  +3.94%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.05'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.19%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.10'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.07%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '155.68'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.48%  '

$ws.Range('E32').Value = '  +1.17%  '

$ws.Range('E33').Value = '  +2.61%  '

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0814'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.04%  '

$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.38'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.65%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.28'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +10.98%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.115'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.68%  '

$ws.Range('E38').Value = '  +1.91%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '23.73'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -13.95%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.43'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.89%  '

$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.89'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.00%  '

$ws.Range('B42').Value = 'ApeXProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.07'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +30.42%  '

$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0310'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.71%  '

$ws.Range('D44').Value = '2.067.78'
$ws.Range('E44').Value = '  +3.23%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.998'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.38%  '

$ws.Range('E46').Value = '  +5.86%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.32'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.35%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '77.25'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +14.38%  '

$ws.Range('D49').Value = '2.822.14'
$ws.Range('E49').Value = '  +3.14%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '106.10'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.89%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.192'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.54%  '

